$wb = $excel.ActiveWorkbook

# --- Update selection on PostFunctionCounter ---
$wsPFC = $wb.Worksheets.Item("PostFunctionCounter")
[void]$wsPFC.Activate()
[void]$wsPFC.Range("A4").Select()

# --- Update selection on ConditionsCounter ---
$wsCC = $wb.Worksheets.Item("ConditionsCounter")
[void]$wsCC.Activate()
[void]$wsCC.Range("A3").Select()

# --- Add the new PostFunctionText worksheet after ValidatorsCounter (the last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PostFunctionText"

# Populate cells in the same order the original authoring used so that
# the shared-strings table lines up with the source workbook.
$ws.Range("A1").Value = "rowXpath"
$ws.Range("B1").Value = "postFunctionsList"
$ws.Range("B2").Value = "//tr[@class='transition-details-row']/td/div/div[@id='glass-transitions-postfunctions-panel-1']/div/div/ol/li"
$ws.Range("A2").Value = "//span[contains(., 'Create')]/ancestor::td[@class='transition-name']"
$ws.Range("A3").Value = "//span[contains(., 'Eating')]/ancestor::td[@class='transition-name']"
$ws.Range("A4").Value = "//span[contains(., 'BarfedBackFood')]/ancestor::td[@class='transition-name']"
$ws.Range("B3").Value = "//tr[@class='transition-details-row']/td/div/div[@id='glass-transitions-postfunctions-panel-2']/div/div/ol/li"
$ws.Range("B4").Value = "//tr[@class='transition-details-row']/td/div/div[@id='glass-transitions-postfunctions-panel-4']/div/div/ol/li"
$ws.Range("C1").Value = "expectedNumber"
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 6

# Column widths matching the source workbook (71, 97, ~16.71 characters)
$ws.Columns.Item(1).ColumnWidth = 70.15
$ws.Columns.Item(2).ColumnWidth = 96.15
$ws.Columns.Item(3).ColumnWidth = 15.83

# Final selection/activation state: PostFunctionText becomes the active tab
[void]$ws.Activate()
[void]$ws.Range("B13").Select()
